# Slide 5, "Content Placeholder 3": the paragraph that lists the
# drawbacks of siloed Digital-Twin solutions is split from one
# level-1 bullet (with inline "→" separators) into four bullets:
#   lvl 1 : "Each solution is independent, ... storage systems:"
#   lvl 2 : "No interoperability between different DTs and their data;"
#   lvl 2 : "limiting the capabilities of DTs"
#   lvl 2 : "need to standardize DTs to facilitate integration"  (bold)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# New paragraph texts (the trailing "need to standardize..." run keeps
# its bold formatting, applied separately below).
$introPara   = "However…"
$lead        = "Each solution is independent, in both data models and storage systems:"
$bulletNoInterop  = "No interoperability between different DTs and their data;"
$bulletLimiting   = "limiting the capabilities of DTs"
$bulletStandardize = "need to standardize DTs to facilitate integration"

# Re-assigning the whole TextRange.Text (with `r` paragraph breaks)
# rebuilds the shape's paragraphs from scratch, so each chunk below
# becomes its own <a:p>.
$tr.Text = $introPara + "`r" + $lead + "`r" + $bulletNoInterop + "`r" + $bulletLimiting + "`r" + $bulletStandardize

# Compute the (1-based) start offset of each new paragraph inside the
# freshly-built text so we can restore indent levels / bold via
# Characters(start, length) sub-ranges.
$startLead        = $introPara.Length + 2
$startNoInterop    = $startLead + $lead.Length + 1
$startLimiting     = $startNoInterop + $bulletNoInterop.Length + 1
$startStandardize  = $startLimiting + $bulletLimiting.Length + 1

# "However…" paragraph stays at its original outline level (0).
$tr.Characters(1, $introPara.Length).IndentLevel = 1

# "Each solution ... storage systems:" stays at level 1 (unchanged).
$tr.Characters($startLead, $lead.Length).IndentLevel = 2

# The three new bullets move to level 2.
$tr.Characters($startNoInterop, $bulletNoInterop.Length).IndentLevel = 3
$tr.Characters($startLimiting, $bulletLimiting.Length).IndentLevel = 3
$tr.Characters($startStandardize, $bulletStandardize.Length).IndentLevel = 3

# Keep the bold emphasis on the final bullet's run.
$tr.Characters($startStandardize, $bulletStandardize.Length).Font.Bold = $true
